$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. First paragraph ("Solar Panel Regression" title): add <w:ilvl val="0"/>
#    inside the existing <w:numPr> (before <w:numId val="0"/>).
#    Word's list-level API is 1-based (level 1 == w:ilvl val="0").
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ListFormat.ListLevelNumber = 1

# ---------------------------------------------------------------------
# 2. Second paragraph (the "Unlock how solar power..." paragraph):
#    move the "_GoBack" bookmark from the middle of the paragraph
#    (just before the "over-fitting" run) to the very end of the
#    paragraph (after the last run, right before the paragraph mark).
#
#    A zero-length Range placed exactly at (paragraph.End - 1) triggers
#    a positioning quirk in this engine, so we work around it by
#    bookmarking a temporary character inserted at the end and then
#    deleting that character - the bookmark correctly collapses to
#    zero length in the right spot once the character is removed.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$p2 = $d.Paragraphs.Item(2)
$endOfText = $p2.Range.End - 1
$tempRange = $d.Range($endOfText, $endOfText)
$tempRange.InsertAfter("Z")

$p2b = $d.Paragraphs.Item(2)
$tempCharPos = $p2b.Range.End - 2
$tempCharRange = $d.Range($tempCharPos, $tempCharPos + 1)
$d.Bookmarks.Add("_GoBack", $tempCharRange)
$tempCharRange2 = $d.Range($tempCharPos, $tempCharPos + 1)
$tempCharRange2.Delete()

# ---------------------------------------------------------------------
# 3. Remove the paragraphs that followed the second paragraph:
#      - one empty paragraph
#      - "Presentation : - <link>"
#      - "GitHub - Script file: - <link>"
#      - "GitHub - Deployment file: - <link>"
#      - "App Link : - <link>"
#      - one trailing empty paragraph
#    These are paragraphs 3 through 8 (paragraph 2 is the body text,
#    paragraph 9 is the document's final empty paragraph which stays).
# ---------------------------------------------------------------------
$firstDoomed = $d.Paragraphs.Item(3)
$lastDoomed = $d.Paragraphs.Item(8)
$killRange = $d.Range($firstDoomed.Range.Start, $lastDoomed.Range.End)
$killRange.Delete()
